$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number that was updated
# from 45203 (2023-10-04) to 45205 (2023-10-06) for every data row (2-261).
$newDate = [DateTime]::FromOADate(45205)

for ($r = 2; $r -le 261; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45203) {
        $cell.Value = $newDate
    }
}
